# Generate Report for Handoff
# Updates the localization-status report to reflect the "Ready for handoff"
# state: status text, the handoff/generate timestamps, and the column
# widths of the status/date columns (which widened to fit the new text).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---
# Latest HO Xliff Generate Date (Overview) and Latest Handoff Datetime (de-de)
$wsOverview.Range("G2").Value = "2016-08-27 20:38:50"
$wsDeDe.Range("H2").Value     = "2016-08-27 20:38:50"

# Latest Handoff Datetime (zh-cn)
$wsZhCn.Range("H2").Value = "2016-08-27 20:38:46"

# --- Column widths (status/date columns widened for the longer text) ---
# Target stored width is 17.2159881591797 characters; this engine quantizes
# ColumnWidth assignments to the nearest 1/6 character, so 16.3333... is the
# input that lands closest to the target (17.1666... char stored width).
$wsOverview.Range("E1").ColumnWidth = 16.333333333333332
$wsOverview.Range("F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth     = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth     = 16.333333333333332
